# Auto-generated edit script applying cryptos.xlsx price/volume updates
# plus two row-pair swaps (Fetch.AI/Binance-PegBSC-USD and InjectiveProtocol/Hedera).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing it to stay TEXT (matches the source
# workbook, which stores every Price/Volume cell as an inline string even
# when the text looks like a plain number, e.g. "0.999" or "20.85").
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

$ws.Range("D2").Value = "63.239.73"
$ws.Range("E2").Value = "  +1.49%  "
$ws.Range("D3").Value = "2.473.79"
$ws.Range("E3").Value = "  +1.75%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.10%  "
Set-TextValue $ws.Range("D5") "574.52"
$ws.Range("E5").Value = "  +1.85%  "
Set-TextValue $ws.Range("D6") "147.19"
$ws.Range("E6").Value = "  +1.68%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.65%  "
$ws.Range("D9").Value = "2.474.33"
$ws.Range("E9").Value = "  +1.80%  "
Set-TextValue $ws.Range("D10") "0.112"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("E11").Value = "  +0.88%  "
Set-TextValue $ws.Range("D12") "5.26"
$ws.Range("E12").Value = "  -0.14%  "
Set-TextValue $ws.Range("D13") "0.357"
$ws.Range("E13").Value = "  +1.63%  "
Set-TextValue $ws.Range("D14") "27.07"
$ws.Range("E14").Value = "  +1.53%  "
Set-TextValue $ws.Range("D15") "0.0000180"
$ws.Range("E15").Value = "  +2.47%  "
$ws.Range("D16").Value = "2.914.94"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").Value = "63.030.28"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").Value = "2.476.37"
$ws.Range("E18").Value = "  +2.19%  "
Set-TextValue $ws.Range("D19") "11.49"
$ws.Range("E19").Value = "  +1.79%  "
Set-TextValue $ws.Range("D20") "7.33"
$ws.Range("E20").Value = "  +6.93%  "
Set-TextValue $ws.Range("D21") "328.84"
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("E22").Value = "  +0.85%  "
Set-TextValue $ws.Range("D23") "1.96"
$ws.Range("E23").Value = "  +12.40%  "
$ws.Range("E24").Value = "  -0.12%  "
Set-TextValue $ws.Range("D25") "65.94"
Set-TextValue $ws.Range("D26") "631.11"
$ws.Range("E26").Value = "  +13.98%  "
$ws.Range("E27").Value = "  +8.74%  "
Set-TextValue $ws.Range("D28") "8.51"
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("E29").Value = "  +1.62%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D30") "1.50"
$ws.Range("E30").Value = "  +5.00%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D31") "0.999"
$ws.Range("E31").Value = "  -0.04%  "
Set-TextValue $ws.Range("D32") "8.27"
$ws.Range("E32").Value = "  -0.25%  "
Set-TextValue $ws.Range("D33") "0.143"
$ws.Range("E33").Value = "  -3.23%  "
$ws.Range("E34").Value = "  +1.72%  "
Set-TextValue $ws.Range("D35") "5.14"
$ws.Range("E35").Value = "  +6.70%  "
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("E37").Value = "  -0.06%  "
Set-TextValue $ws.Range("D38") "0.384"
$ws.Range("E38").Value = "  +0.48%  "
Set-TextValue $ws.Range("D39") "5.43"
$ws.Range("E39").Value = "  -2.80%  "
Set-TextValue $ws.Range("D40") "18.79"
$ws.Range("E40").Value = "  +0.39%  "
Set-TextValue $ws.Range("D41") "147.31"
$ws.Range("E41").Value = "  -1.83%  "
Set-TextValue $ws.Range("D42") "1.79"
$ws.Range("E42").Value = "  -0.52%  "
Set-TextValue $ws.Range("D43") "2.62"
$ws.Range("E43").Value = "  +12.93%  "
$ws.Range("E44").Value = "  +0.04%  "
Set-TextValue $ws.Range("D45") "148.13"
$ws.Range("E45").Value = "  +0.13%  "
Set-TextValue $ws.Range("D46") "3.73"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D47") "20.85"
$ws.Range("E47").Value = "  +2.78%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D48") "0.0540"
$ws.Range("E48").Value = "  +1.10%  "
Set-TextValue $ws.Range("D49") "0.607"
$ws.Range("E49").Value = "  +1.62%  "
Set-TextValue $ws.Range("D50") "0.0234"
$ws.Range("E50").Value = "  +1.60%  "
Set-TextValue $ws.Range("D51") "0.0923"
$ws.Range("E51").Value = "  -0.27%  "
